$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scotland")
$ws.Range("A2:B2").EntireRow.Insert()
$row = $ws.Rows.Item(2)
$row.Style = "Normal"
$row.Font.Size = 14

$r = $ws.Range("A2:B2")
$r.Font.Bold = $true
